$wb = $excel.ActiveWorkbook

# ---- Matches_SOG: append 4 new match rows (443-446) ----
$wsMatches = $wb.Worksheets.Item("Matches_SOG")

# Row 443
$wsMatches.Range("A443").NumberFormat = "@"
$wsMatches.Range("A443").Value = "897740"
$wsMatches.Range("A443").Style = "Normal"
$wsMatches.Range("B443").Value = "2025-11-07T12:15:00"
$wsMatches.Range("C443").Value = "Амур"
$wsMatches.Range("D443").Value = "Металлург Мг"
$wsMatches.Range("E443").Value = 29
$wsMatches.Range("F443").Value = 31
$wsMatches.Range("G443").Value = "khl_text"

# Row 444
$wsMatches.Range("A444").NumberFormat = "@"
$wsMatches.Range("A444").Value = "897743"
$wsMatches.Range("A444").Style = "Normal"
$wsMatches.Range("B444").Value = "2025-11-07T12:30:00"
$wsMatches.Range("C444").Value = "Адмирал"
$wsMatches.Range("D444").Value = "Авангард"
$wsMatches.Range("E444").Value = 37
$wsMatches.Range("F444").Value = 24
$wsMatches.Range("G444").Value = "khl_text"

# Row 445
$wsMatches.Range("A445").NumberFormat = "@"
$wsMatches.Range("A445").Value = "897742"
$wsMatches.Range("A445").Style = "Normal"
$wsMatches.Range("B445").Value = "2025-11-07T17:30:00"
$wsMatches.Range("C445").Value = "Барыс"
$wsMatches.Range("D445").Value = "Автомобилист"
$wsMatches.Range("E445").Value = 40
$wsMatches.Range("F445").Value = 24
$wsMatches.Range("G445").Value = "khl_text"

# Row 446
$wsMatches.Range("A446").NumberFormat = "@"
$wsMatches.Range("A446").Value = "897741"
$wsMatches.Range("A446").Style = "Normal"
$wsMatches.Range("B446").Value = "2025-11-07T19:00:00"
$wsMatches.Range("C446").Value = "Нефтехимик"
$wsMatches.Range("D446").Value = "Салават Юлаев"
$wsMatches.Range("E446").Value = 25
$wsMatches.Range("F446").Value = 36
$wsMatches.Range("G446").Value = "khl_text"

# ---- Shots_HA: refresh as_of_utc + aggregate stats ----
$wsShotsHA = $wb.Worksheets.Item("Shots_HA")

# Row 2
$wsShotsHA.Range("D2").Value = "2025-11-07T19:00:00Z"
$wsShotsHA.Range("F2").Value = 18
$wsShotsHA.Range("K2").Value = 642
$wsShotsHA.Range("L2").Value = 550
$wsShotsHA.Range("M2").Value = 35.7
$wsShotsHA.Range("N2").Value = 30.6

# Row 3
$wsShotsHA.Range("D3").Value = "2025-11-07T19:00:00Z"
$wsShotsHA.Range("F3").Value = 25
$wsShotsHA.Range("K3").Value = 692
$wsShotsHA.Range("L3").Value = 765
$wsShotsHA.Range("M3").Value = 27.7
$wsShotsHA.Range("N3").Value = 30.6

# Row 4
$wsShotsHA.Range("D4").Value = "2025-11-07T19:00:00Z"
$wsShotsHA.Range("E4").Value = 16
$wsShotsHA.Range("G4").Value = 617
$wsShotsHA.Range("H4").Value = 436
$wsShotsHA.Range("I4").Value = 38.6
$wsShotsHA.Range("J4").Value = 27.2

# Row 5
$wsShotsHA.Range("D5").Value = "2025-11-07T19:00:00Z"

# Row 6
$wsShotsHA.Range("D6").Value = "2025-11-07T19:00:00Z"
$wsShotsHA.Range("E6").Value = 19
$wsShotsHA.Range("G6").Value = 581
$wsShotsHA.Range("H6").Value = 671
$wsShotsHA.Range("I6").Value = 30.6
$wsShotsHA.Range("J6").Value = 35.3

# Row 7
$wsShotsHA.Range("D7").Value = "2025-11-07T19:00:00Z"
$wsShotsHA.Range("E7").Value = 29
$wsShotsHA.Range("G7").Value = 925
$wsShotsHA.Range("H7").Value = 919
$wsShotsHA.Range("I7").Value = 31.9
$wsShotsHA.Range("J7").Value = 31.7

# Row 8
$wsShotsHA.Range("D8").Value = "2025-11-07T19:00:00Z"

# Row 9
$wsShotsHA.Range("D9").Value = "2025-11-07T19:00:00Z"

# Row 10
$wsShotsHA.Range("D10").Value = "2025-11-07T19:00:00Z"

# Row 11
$wsShotsHA.Range("D11").Value = "2025-11-07T19:00:00Z"

# Row 12
$wsShotsHA.Range("D12").Value = "2025-11-07T19:00:00Z"

# Row 13
$wsShotsHA.Range("D13").Value = "2025-11-07T19:00:00Z"
$wsShotsHA.Range("F13").Value = 16
$wsShotsHA.Range("K13").Value = 454
$wsShotsHA.Range("L13").Value = 428
$wsShotsHA.Range("M13").Value = 28.4
$wsShotsHA.Range("N13").Value = 26.8

# Row 14
$wsShotsHA.Range("D14").Value = "2025-11-07T19:00:00Z"
$wsShotsHA.Range("E14").Value = 24
$wsShotsHA.Range("G14").Value = 749
$wsShotsHA.Range("H14").Value = 806
$wsShotsHA.Range("I14").Value = 31.2
$wsShotsHA.Range("J14").Value = 33.6

# Row 15
$wsShotsHA.Range("D15").Value = "2025-11-07T19:00:00Z"

# Row 16
$wsShotsHA.Range("D16").Value = "2025-11-07T19:00:00Z"
$wsShotsHA.Range("F16").Value = 26
$wsShotsHA.Range("K16").Value = 730
$wsShotsHA.Range("L16").Value = 758
$wsShotsHA.Range("M16").Value = 28.1
$wsShotsHA.Range("N16").Value = 29.2

# Row 17
$wsShotsHA.Range("D17").Value = "2025-11-07T19:00:00Z"

# Row 18
$wsShotsHA.Range("D18").Value = "2025-11-07T19:00:00Z"

# Row 19
$wsShotsHA.Range("D19").Value = "2025-11-07T19:00:00Z"

# Row 20
$wsShotsHA.Range("D20").Value = "2025-11-07T19:00:00Z"

# Row 21
$wsShotsHA.Range("D21").Value = "2025-11-07T19:00:00Z"

# Row 22
$wsShotsHA.Range("D22").Value = "2025-11-07T19:00:00Z"

# Row 23
$wsShotsHA.Range("D23").Value = "2025-11-07T19:00:00Z"

# ---- Shots_Summary: refresh as_of_utc + aggregate stats ----
$wsShotsSummary = $wb.Worksheets.Item("Shots_Summary")

# Row 2
$wsShotsSummary.Range("D2").Value = "2025-11-07T19:00:00Z"
$wsShotsSummary.Range("E2").Value = 39
$wsShotsSummary.Range("F2").Value = 1329
$wsShotsSummary.Range("G2").Value = 1140
$wsShotsSummary.Range("H2").Value = 34.1
$wsShotsSummary.Range("I2").Value = 29.2

# Row 3
$wsShotsSummary.Range("D3").Value = "2025-11-07T19:00:00Z"
$wsShotsSummary.Range("E3").Value = 43
$wsShotsSummary.Range("F3").Value = 1213
$wsShotsSummary.Range("G3").Value = 1323
$wsShotsSummary.Range("H3").Value = 28.2
$wsShotsSummary.Range("I3").Value = 30.8

# Row 4
$wsShotsSummary.Range("D4").Value = "2025-11-07T19:00:00Z"
$wsShotsSummary.Range("E4").Value = 36
$wsShotsSummary.Range("F4").Value = 1255
$wsShotsSummary.Range("G4").Value = 997
$wsShotsSummary.Range("H4").Value = 34.9
$wsShotsSummary.Range("I4").Value = 27.7

# Row 5
$wsShotsSummary.Range("D5").Value = "2025-11-07T19:00:00Z"

# Row 6
$wsShotsSummary.Range("D6").Value = "2025-11-07T19:00:00Z"
$wsShotsSummary.Range("E6").Value = 40
$wsShotsSummary.Range("F6").Value = 1173
$wsShotsSummary.Range("G6").Value = 1440
$wsShotsSummary.Range("I6").Value = 36

# Row 7
$wsShotsSummary.Range("D7").Value = "2025-11-07T19:00:00Z"
$wsShotsSummary.Range("E7").Value = 44
$wsShotsSummary.Range("F7").Value = 1340
$wsShotsSummary.Range("G7").Value = 1434
$wsShotsSummary.Range("H7").Value = 30.5
$wsShotsSummary.Range("I7").Value = 32.6

# Row 8
$wsShotsSummary.Range("D8").Value = "2025-11-07T19:00:00Z"

# Row 9
$wsShotsSummary.Range("D9").Value = "2025-11-07T19:00:00Z"

# Row 10
$wsShotsSummary.Range("D10").Value = "2025-11-07T19:00:00Z"

# Row 11
$wsShotsSummary.Range("D11").Value = "2025-11-07T19:00:00Z"

# Row 12
$wsShotsSummary.Range("D12").Value = "2025-11-07T19:00:00Z"

# Row 13
$wsShotsSummary.Range("D13").Value = "2025-11-07T19:00:00Z"
$wsShotsSummary.Range("E13").Value = 41
$wsShotsSummary.Range("F13").Value = 1338
$wsShotsSummary.Range("G13").Value = 1074
$wsShotsSummary.Range("H13").Value = 32.6
$wsShotsSummary.Range("I13").Value = 26.2

# Row 14
$wsShotsSummary.Range("D14").Value = "2025-11-07T19:00:00Z"
$wsShotsSummary.Range("E14").Value = 43
$wsShotsSummary.Range("F14").Value = 1266
$wsShotsSummary.Range("G14").Value = 1528
$wsShotsSummary.Range("H14").Value = 29.4

# Row 15
$wsShotsSummary.Range("D15").Value = "2025-11-07T19:00:00Z"

# Row 16
$wsShotsSummary.Range("D16").Value = "2025-11-07T19:00:00Z"
$wsShotsSummary.Range("E16").Value = 41
$wsShotsSummary.Range("F16").Value = 1133
$wsShotsSummary.Range("G16").Value = 1176
$wsShotsSummary.Range("H16").Value = 27.6
$wsShotsSummary.Range("I16").Value = 28.7

# Row 17
$wsShotsSummary.Range("D17").Value = "2025-11-07T19:00:00Z"

# Row 18
$wsShotsSummary.Range("D18").Value = "2025-11-07T19:00:00Z"

# Row 19
$wsShotsSummary.Range("D19").Value = "2025-11-07T19:00:00Z"

# Row 20
$wsShotsSummary.Range("D20").Value = "2025-11-07T19:00:00Z"

# Row 21
$wsShotsSummary.Range("D21").Value = "2025-11-07T19:00:00Z"

# Row 22
$wsShotsSummary.Range("D22").Value = "2025-11-07T19:00:00Z"

# Row 23
$wsShotsSummary.Range("D23").Value = "2025-11-07T19:00:00Z"

# ---- Meta_ext: bump as_of_utc + build_version ----
$wsMeta = $wb.Worksheets.Item("Meta_ext")
$wsMeta.Range("B2").Value = "2025-11-07T19:00:00Z"
$wsMeta.Range("D2").Value = 52
